$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stat 1 table (rows 25-33): add TL / TR columns (K, L) ---

# Header row 25: labels "TL" / "TR" (existing cells already carry style s="3")
$ws.Range("K25").Value = "TL"
$ws.Range("L25").Value = "TR"

# Header row 26: unit "[N]" for both (existing cells already carry style s="3")
$ws.Range("K26").Value = "[N]"
$ws.Range("L26").Value = "[N]"

# Row 27 (blank separator row): K27/L27 cells are removed entirely
$ws.Range("K27").Clear()
$ws.Range("L27").Clear()

# Data rows 28-33: new numeric values, no explicit style (default)
$ws.Range("K28").Clear()
$ws.Range("K28").Value = 3106.04
$ws.Range("L28").Clear()
$ws.Range("L28").Value = 3413.14

$ws.Range("K29").Clear()
$ws.Range("K29").Value = 3086.32
$ws.Range("L29").Clear()
$ws.Range("L29").Value = 3361.08

$ws.Range("K30").Clear()
$ws.Range("K30").Value = 2323.0500000000002
$ws.Range("L30").Clear()
$ws.Range("L30").Value = 2616.9499999999998

$ws.Range("K31").Clear()
$ws.Range("K31").Value = 1811.51
$ws.Range("L31").Clear()
$ws.Range("L31").Value = 2174.62

$ws.Range("K32").Clear()
$ws.Range("K32").Value = 1714.73
$ws.Range("L32").Clear()
$ws.Range("L32").Value = 1981.96

$ws.Range("K33").Clear()
$ws.Range("K33").Value = 2088.69
$ws.Range("L33").Clear()
$ws.Range("L33").Value = 2469.9299999999998

# --- Stat 2 table (rows 56-63): add TL / TR / Tps1eng columns (N, O, P) ---

# Header row 56: copy formatting from existing M56 cell (style s="3"), then set labels
$ws.Range("M56").Copy()
$ws.Range("N56").PasteSpecial(-4122)
$ws.Range("O56").PasteSpecial(-4122)
$ws.Range("P56").PasteSpecial(-4122)
$ws.Range("N56").Value = "TL"
$ws.Range("O56").Value = "TR"
$ws.Range("P56").Value = "Tps1eng"

# Header row 57: same style s="3" (copied from M56, not M57), units "[N]"
$ws.Range("M56").Copy()
$ws.Range("N57").PasteSpecial(-4122)
$ws.Range("O57").PasteSpecial(-4122)
$ws.Range("P57").PasteSpecial(-4122)
$ws.Range("N57").Value = "[N]"
$ws.Range("O57").Value = "[N]"
$ws.Range("P57").Value = "[N]"

# Data rows 59-63: new numeric values, no explicit style (default)
$ws.Range("N59").Value = 1947.92
$ws.Range("O59").Value = 2242.0700000000002
$ws.Range("P59").Value = 1427.77

$ws.Range("N60").Value = 1979.8
$ws.Range("O60").Value = 2278.42
$ws.Range("P60").Value = 1479.82

$ws.Range("N61").Value = 2016.27
$ws.Range("O61").Value = 2327.48
$ws.Range("P61").Value = 1537.34

$ws.Range("N62").Value = 1931.49
$ws.Range("O62").Value = 2207.48
$ws.Range("P62").Value = 1360.7

$ws.Range("N63").Value = 1917.88
$ws.Range("O63").Value = 2196.9299999999998
$ws.Range("P63").Value = 1308.32

# --- View state: update the saved selection to match the author's last position ---
$ws.Range("L54").Select()
